$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Hoja1" to "Datos"
$ws.Name = "Datos"

# Change header cell A1 from "Caso de Prueba" to "TestCase"
$ws.Range("A1").Value = "TestCase"

# Reset the cursor/selection back to the top-left cell (the original file had
# a stray "F3" selection left over from editing; move it back to A1).
$ws.Range("A1").Select() | Out-Null
